$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document. Word.InsertParagraphAfter
#    would inherit the Heading1 style of the title, so the new paragraph's
#    content (including its lack of an explicit paragraph style, a leading
#    empty run, a bold "Meta description" run, and a plain run with the rest
#    of the sentence) is supplied directly via InsertXML.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs(2)

$metaXml = '<?xml version="1.0" standalone="yes"?>' +
  '<?mso-application progid="Word.Document"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' +
              '<w:r/>' +
              '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +
              '<w:r><w:t>: Experience the serene atmosphere of Buffalo Spirit Slot with two Wild symbols and a Bonus mode. Try it for free now!</w:t></w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

[void]$metaPara.Range.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# 2) Remove the duplicated bold title paragraph that used to sit near the end
#    of the document (directly before the italic "meta description -> image
#    prompt" paragraph). After step 1 inserted exactly one new paragraph,
#    this duplicate is always the second-to-last paragraph in the document.
# ---------------------------------------------------------------------------
$dupIndex = $d.Paragraphs.Count - 1
$dupPara = $d.Paragraphs($dupIndex)
$dupPara.Range.Delete()

# ---------------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph with the new image
#    generation prompt, preserving its run formatting (and leading empty
#    run). The "- 1" keeps the trailing paragraph mark out of the replaced
#    range so the paragraph count/formatting is untouched.
# ---------------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
$newRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$newRange.Text = "Create a unique feature image for Buffalo Spirit that captures the calmness and serenity of the game while also featuring a happy Maya warrior with glasses. The image should be in a cartoon style. The warrior should convey a sense of joy and relaxation while playing the game. The background should include elements that evoke the American West, such as cactus plants and desert rocks. The buffalo should also feature prominently in the image, with one of the Wild symbols being highlighted. The overall vibe of the image should be delightful and inviting, encouraging potential players to engage in a stress-free gaming experience with Buffalo Spirit."
